$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1560.2396
$ws.Range("I15").Value = 1560.2396
$ws.Range("K15").Value = 4680.718800000001
$ws.Range("M15").Value = -4511.718800000001

$ws.Range("H138").Value = 3159.361
$ws.Range("I138").Value = 1384.5676
$ws.Range("J138").Value = 5035.5713
$ws.Range("K138").Value = 4153.7028
$ws.Range("L138").Value = 15106.7139
$ws.Range("M138").Value = 986.2972
$ws.Range("N138").Value = -25386.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13480.124
$ws.Range("I32").Value = 14184.047
$ws.Range("J32").Value = 7976.727
$ws.Range("K32").Value = 14184.047
$ws.Range("L32").Value = 7976.727
$ws.Range("M32").Value = -13897.047
$ws.Range("N32").Value = -8550.726999999999

$ws.Range("H61").Value = 5434.948
$ws.Range("I61").Value = 3046.1875
$ws.Range("J61").Value = 17195
$ws.Range("K61").Value = 3046.1875
$ws.Range("L61").Value = 17195
$ws.Range("M61").Value = -2834.1875
$ws.Range("N61").Value = -17619

$ws.Range("H86").Value = 38412
$ws.Range("J86").Value = 38412
$ws.Range("L86").Value = 38412
$ws.Range("N86").Value = -40784

$ws.Range("H89").Value = 38412
$ws.Range("J89").Value = 38412
$ws.Range("L89").Value = 115236
$ws.Range("N89").Value = -127092

$ws.Range("H136").Value = 5434.948
$ws.Range("I136").Value = 3046.1875
$ws.Range("J136").Value = 17195
$ws.Range("K136").Value = 9138.5625
$ws.Range("L136").Value = 51585
$ws.Range("M136").Value = -6588.5625
$ws.Range("N136").Value = -56685

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 450
$ws.Range("I22").Value = 425
$ws.Range("K22").Value = 425
$ws.Range("M22").Value = -252

$ws.Range("H134").Value = 22308.64
$ws.Range("I134").Value = 2236.8572
$ws.Range("J134").Value = 69142.8
$ws.Range("K134").Value = 6710.571599999999
$ws.Range("L134").Value = 207428.4
$ws.Range("M134").Value = -4175.571599999999
$ws.Range("N134").Value = -212498.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1110072.1
$ws.Range("I58").Value = 1491139.9
$ws.Range("J58").Value = 3161.1904
$ws.Range("K58").Value = 1491139.9
$ws.Range("L58").Value = 3161.1904
$ws.Range("M58").Value = -1490936.9
$ws.Range("N58").Value = -3567.1904

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 3634.3635
$ws.Range("I132").Value = 3812.5264
$ws.Range("J132").Value = 3236.1177
$ws.Range("K132").Value = 11437.5792
$ws.Range("L132").Value = 9708.3531
$ws.Range("M132").Value = -8907.5792
$ws.Range("N132").Value = -14768.3531

$ws.Range("H134").Value = 2115.625
$ws.Range("I134").Value = 1301.174
$ws.Range("J134").Value = 3217.5293
$ws.Range("K134").Value = 3903.522
$ws.Range("L134").Value = 9652.5879
$ws.Range("M134").Value = -1368.522
$ws.Range("N134").Value = -14722.5879

$ws.Range("H136").Value = 1110072.1
$ws.Range("I136").Value = 1491139.9
$ws.Range("J136").Value = 3161.1904
$ws.Range("K136").Value = 4473419.699999999
$ws.Range("L136").Value = 9483.5712
$ws.Range("M136").Value = -4470869.699999999
$ws.Range("N136").Value = -14583.5712

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7247219
$ws.Range("I5").Value = 441.13333
$ws.Range("J5").Value = 20834928
$ws.Range("K5").Value = 1323.39999
$ws.Range("L5").Value = 62504784
$ws.Range("M5").Value = -1211.39999
$ws.Range("N5").Value = -62505008

$ws.Range("H50").Value = 260.10526
$ws.Range("I50").Value = 150.28572
$ws.Range("J50").Value = 324.16666
$ws.Range("K50").Value = 450.85716
$ws.Range("L50").Value = 972.4999799999999
$ws.Range("M50").Value = 30.14283999999998
$ws.Range("N50").Value = -1934.49998

$ws.Range("H53").Value = 260.10526
$ws.Range("I53").Value = 150.28572
$ws.Range("J53").Value = 324.16666
$ws.Range("K53").Value = 450.85716
$ws.Range("L53").Value = 972.4999799999999
$ws.Range("M53").Value = 30.14283999999998
$ws.Range("N53").Value = -1934.49998

$ws.Range("H58").Value = 3100
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 3100
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 9300
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -9556

$ws.Range("H122").Value = 825.9643
$ws.Range("I122").Value = 576.1667
$ws.Range("J122").Value = 1013.3125
$ws.Range("K122").Value = 5185.5003
$ws.Range("L122").Value = 9119.8125
$ws.Range("M122").Value = -2735.5003
$ws.Range("N122").Value = -14019.8125

$ws.Range("H135").Value = 7247219
$ws.Range("I135").Value = 441.13333
$ws.Range("J135").Value = 20834928
$ws.Range("K135").Value = 3970.19997
$ws.Range("L135").Value = 187514352
$ws.Range("M135").Value = -1435.19997
$ws.Range("N135").Value = -187519422

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23625
$ws.Range("I46").Value = 7500
$ws.Range("J46").Value = 29000
$ws.Range("K46").Value = 7500
$ws.Range("L46").Value = 29000
$ws.Range("M46").Value = -7344
$ws.Range("N46").Value = -29312

$ws.Range("H132").Value = 3488.8086
$ws.Range("I132").Value = 1487.2
$ws.Range("J132").Value = 9326.833000000001
$ws.Range("K132").Value = 4461.6
$ws.Range("L132").Value = 27980.499
$ws.Range("M132").Value = -1931.6
$ws.Range("N132").Value = -33040.499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 9800
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H46").Value = 725.125
$ws.Range("I46").Value = 400.25
$ws.Range("J46").Value = 1050
$ws.Range("K46").Value = 400.25
$ws.Range("L46").Value = 1050
$ws.Range("M46").Value = -212.25
$ws.Range("N46").Value = -1426

$ws.Range("H122").Value = 7332.135
$ws.Range("I122").Value = 7011.914
$ws.Range("J122").Value = 7991.4116
$ws.Range("K122").Value = 21035.742
$ws.Range("L122").Value = 23974.2348
$ws.Range("M122").Value = -18585.742
$ws.Range("N122").Value = -28874.2348

$ws.Range("H136").Value = 3214.2727
$ws.Range("I136").Value = 1812.6552
$ws.Range("J136").Value = 7492.8945
$ws.Range("K136").Value = 5437.9656
$ws.Range("L136").Value = 22478.6835
$ws.Range("M136").Value = -2887.9656
$ws.Range("N136").Value = -27578.6835

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2510.276
$ws.Range("I122").Value = 1879.72
$ws.Range("J122").Value = 6451.25
$ws.Range("K122").Value = 5639.16
$ws.Range("L122").Value = 19353.75
$ws.Range("M122").Value = -3189.16
$ws.Range("N122").Value = -24253.75

$ws.Range("H132").Value = 1229.2727
$ws.Range("I132").Value = 644.1316
$ws.Range("J132").Value = 2537.2354
$ws.Range("K132").Value = 1932.3948
$ws.Range("L132").Value = 7611.706200000001
$ws.Range("M132").Value = 597.6052
$ws.Range("N132").Value = -12671.7062
